$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete old data row 84 (oldest entry, "M.15608.6479.683840").
# This shifts the summary row (old row 85) up to become row 84, preserving its style,
# and removes the now-obsolete last data row in one operation.
$ws.Rows.Item(84).Delete()

# Step 2: write the 8 brand-new rows (2-9) and the shifted historical rows (10-83).
# Column I (SharedId) is always blank across the whole sheet, so it is left untouched
# (it already holds an empty value with the correct per-row style).
$ws.Range("A2").Value = 45935.92592592593
$ws.Range("B2").Value = "M.15608.6487.568983"
$ws.Range("C2").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D2").Value = "Sale (iOS App)"
$ws.Range("E2").Value = "Pending"
$ws.Range("F2").Value = 88.39
$ws.Range("G2").Value = 8.84
$ws.Range("H2").Value = "SAHS"

$ws.Range("A3").Value = 45935.850381944445
$ws.Range("B3").Value = "M.15608.6487.461322"
$ws.Range("C3").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D3").Value = "Sale (iOS App)"
$ws.Range("E3").Value = "Pending"
$ws.Range("F3").Value = 29.99
$ws.Range("G3").Value = 3.0
$ws.Range("H3").Value = "SAHS"

$ws.Range("A4").Value = 45935.677824074075
$ws.Range("B4").Value = "M.15608.6487.95086"
$ws.Range("C4").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D4").Value = "Sale (iOS App)"
$ws.Range("E4").Value = "Pending"
$ws.Range("F4").Value = 28.41
$ws.Range("G4").Value = 2.84
$ws.Range("H4").Value = "SAHS"

$ws.Range("A5").Value = 45935.59841435185
$ws.Range("B5").Value = "M.15608.6487.44991"
$ws.Range("C5").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D5").Value = "Sale (iOS App)"
$ws.Range("E5").Value = "Pending"
$ws.Range("F5").Value = 26.83
$ws.Range("G5").Value = 2.68
$ws.Range("H5").Value = "SAHS15"

$ws.Range("A6").Value = 45935.55547453704
$ws.Range("B6").Value = "M.15608.6487.26758"
$ws.Range("C6").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D6").Value = "Sale (iOS App)"
$ws.Range("E6").Value = "Pending"
$ws.Range("F6").Value = 86.82
$ws.Range("G6").Value = 8.68
$ws.Range("H6").Value = "SAHS"

$ws.Range("A7").Value = 45935.541284722225
$ws.Range("B7").Value = "M.15608.6487.21349"
$ws.Range("C7").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D7").Value = "Sale (Android App)"
$ws.Range("E7").Value = "Pending"
$ws.Range("F7").Value = 39.46
$ws.Range("G7").Value = 3.95
$ws.Range("H7").Value = "SAHS"

$ws.Range("A8").Value = 45935.3608912037
$ws.Range("B8").Value = "M.15608.6486.1199129"
$ws.Range("C8").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D8").Value = "Sale (iOS App)"
$ws.Range("E8").Value = "Pending"
$ws.Range("F8").Value = 64.72
$ws.Range("G8").Value = 6.47
$ws.Range("H8").Value = "SAHS"

$ws.Range("A9").Value = 45935.357453703706
$ws.Range("B9").Value = "M.15608.6486.1184171"
$ws.Range("C9").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D9").Value = "Sale (iOS App)"
$ws.Range("E9").Value = "Pending"
$ws.Range("F9").Value = 55.14
$ws.Range("G9").Value = 5.51
$ws.Range("H9").Value = "SAHS15"

$ws.Range("A10").Value = 45934.90261574074
$ws.Range("B10").Value = "M.15608.6486.553387"
$ws.Range("C10").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D10").Value = "Sale (iOS App)"
$ws.Range("E10").Value = "Pending"
$ws.Range("F10").Value = 47.35
$ws.Range("G10").Value = 4.74
$ws.Range("H10").Value = "SAHS"

$ws.Range("A11").Value = 45934.85689814815
$ws.Range("B11").Value = "M.15608.6486.472769"
$ws.Range("C11").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D11").Value = "Sale (iOS App)"
$ws.Range("E11").Value = "Pending"
$ws.Range("F11").Value = 36.3
$ws.Range("G11").Value = 3.63
$ws.Range("H11").Value = "SAHS"

$ws.Range("A12").Value = 45934.804236111115
$ws.Range("B12").Value = "M.15608.6486.377939"
$ws.Range("C12").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D12").Value = "Sale (iOS App)"
$ws.Range("E12").Value = "Pending"
$ws.Range("F12").Value = 86.82
$ws.Range("G12").Value = 8.68
$ws.Range("H12").Value = "SAHS"

$ws.Range("A13").Value = 45934.66096064815
$ws.Range("B13").Value = "M.15608.6486.161300"
$ws.Range("C13").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D13").Value = "Sale (iOS App)"
$ws.Range("E13").Value = "Pending"
$ws.Range("F13").Value = 112.07
$ws.Range("G13").Value = 11.21
$ws.Range("H13").Value = "SAHS"

$ws.Range("A14").Value = 45934.63767361111
$ws.Range("B14").Value = "M.15608.6486.135911"
$ws.Range("C14").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D14").Value = "Sale (iOS App)"
$ws.Range("E14").Value = "Pending"
$ws.Range("F14").Value = 12.63
$ws.Range("G14").Value = 1.26
$ws.Range("H14").Value = "SAHS"

$ws.Range("A15").Value = 45934.596979166665
$ws.Range("B15").Value = "M.15608.6486.94009"
$ws.Range("C15").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D15").Value = "Sale (iOS App)"
$ws.Range("E15").Value = "Pending"
$ws.Range("F15").Value = 29.99
$ws.Range("G15").Value = 3.0
$ws.Range("H15").Value = "SAHS"

$ws.Range("A16").Value = 45934.582453703704
$ws.Range("B16").Value = "M.15608.6486.78318"
$ws.Range("C16").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D16").Value = "Sale (iOS App)"
$ws.Range("E16").Value = "Pending"
$ws.Range("F16").Value = 14.21
$ws.Range("G16").Value = 1.42
$ws.Range("H16").Value = "SAHS"

$ws.Range("A17").Value = 45934.545
$ws.Range("B17").Value = "M.15608.6486.39475"
$ws.Range("C17").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D17").Value = "Sale (iOS App)"
$ws.Range("E17").Value = "Pending"
$ws.Range("F17").Value = 14.21
$ws.Range("G17").Value = 1.42
$ws.Range("H17").Value = "SAHS"

$ws.Range("A18").Value = 45934.53135416667
$ws.Range("B18").Value = "M.15608.6486.30882"
$ws.Range("C18").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D18").Value = "Sale (iOS App)"
$ws.Range("E18").Value = "Pending"
$ws.Range("F18").Value = 65.59
$ws.Range("G18").Value = 6.56
$ws.Range("H18").Value = "SAHS15"

$ws.Range("A19").Value = 45934.503125
$ws.Range("B19").Value = "M.15608.6486.3193"
$ws.Range("C19").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D19").Value = "Sale (iOS App)"
$ws.Range("E19").Value = "Pending"
$ws.Range("F19").Value = 14.21
$ws.Range("G19").Value = 1.42
$ws.Range("H19").Value = "SAHS"

$ws.Range("A20").Value = 45934.18512731481
$ws.Range("B20").Value = "M.15608.6485.1063730"
$ws.Range("C20").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D20").Value = "Sale (iOS App)"
$ws.Range("E20").Value = "Pending"
$ws.Range("F20").Value = 83.48
$ws.Range("G20").Value = 8.35
$ws.Range("H20").Value = "SAHS15"

$ws.Range("A21").Value = 45934.05074074074
$ws.Range("B21").Value = "M.15608.6485.784899"
$ws.Range("C21").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D21").Value = "Sale (iOS App)"
$ws.Range("E21").Value = "Pending"
$ws.Range("F21").Value = 14.21
$ws.Range("G21").Value = 1.42
$ws.Range("H21").Value = "SAHS"

$ws.Range("A22").Value = 45934.02245370371
$ws.Range("B22").Value = "M.15608.6485.729554"
$ws.Range("C22").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D22").Value = "Sale (iOS App)"
$ws.Range("E22").Value = "Pending"
$ws.Range("F22").Value = 179.95
$ws.Range("G22").Value = 17.99
$ws.Range("H22").Value = "SAHS"

$ws.Range("A23").Value = 45934.02226851852
$ws.Range("B23").Value = "M.15608.6485.729758"
$ws.Range("C23").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D23").Value = "Sale (iOS App)"
$ws.Range("E23").Value = "Pending"
$ws.Range("F23").Value = 110.49
$ws.Range("G23").Value = 11.05
$ws.Range("H23").Value = "SAHS"

$ws.Range("A24").Value = 45934.00570601852
$ws.Range("B24").Value = "M.15608.6485.708442"
$ws.Range("C24").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D24").Value = "Sale (iOS App)"
$ws.Range("E24").Value = "Pending"
$ws.Range("F24").Value = 110.49
$ws.Range("G24").Value = 11.05
$ws.Range("H24").Value = "SAHS"

$ws.Range("A25").Value = 45933.99395833333
$ws.Range("B25").Value = "M.15608.6485.672663"
$ws.Range("C25").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D25").Value = "Sale (iOS App)"
$ws.Range("E25").Value = "Pending"
$ws.Range("F25").Value = 71.03
$ws.Range("G25").Value = 7.1
$ws.Range("H25").Value = "SAHS"

$ws.Range("A26").Value = 45933.98636574074
$ws.Range("B26").Value = "M.15608.6485.663296"
$ws.Range("C26").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D26").Value = "Sale (Android App)"
$ws.Range("E26").Value = "Pending"
$ws.Range("F26").Value = 61.56
$ws.Range("G26").Value = 6.16
$ws.Range("H26").Value = "SAHS"

$ws.Range("A27").Value = 45933.88585648148
$ws.Range("B27").Value = "M.15608.6485.463996"
$ws.Range("C27").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D27").Value = "Sale (iOS App)"
$ws.Range("E27").Value = "Pending"
$ws.Range("F27").Value = 14.21
$ws.Range("G27").Value = 1.42
$ws.Range("H27").Value = "SAHS"

$ws.Range("A28").Value = 45933.81726851852
$ws.Range("B28").Value = "M.15608.6485.350281"
$ws.Range("C28").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D28").Value = "Sale (iOS App)"
$ws.Range("E28").Value = "Pending"
$ws.Range("F28").Value = 93.13
$ws.Range("G28").Value = 9.31
$ws.Range("H28").Value = "SAHS"

$ws.Range("A29").Value = 45933.75896990741
$ws.Range("B29").Value = "M.15608.6485.234639"
$ws.Range("C29").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D29").Value = "Sale (iOS App)"
$ws.Range("E29").Value = "Pending"
$ws.Range("F29").Value = 14.21
$ws.Range("G29").Value = 1.42
$ws.Range("H29").Value = "SAHS"

$ws.Range("A30").Value = 45933.75038194445
$ws.Range("B30").Value = "M.15608.6485.224423"
$ws.Range("C30").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D30").Value = "Sale (iOS App)"
$ws.Range("E30").Value = "Pending"
$ws.Range("F30").Value = 52.16
$ws.Range("G30").Value = 5.22
$ws.Range("H30").Value = "SAHS15"

$ws.Range("A31").Value = 45933.748252314814
$ws.Range("B31").Value = "M.15608.6485.222948"
$ws.Range("C31").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D31").Value = "Sale (iOS App)"
$ws.Range("E31").Value = "Pending"
$ws.Range("F31").Value = 26.83
$ws.Range("G31").Value = 2.68
$ws.Range("H31").Value = "SAHS"

$ws.Range("A32").Value = 45933.60122685185
$ws.Range("B32").Value = "M.15608.6485.53178"
$ws.Range("C32").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D32").Value = "Sale (iOS App)"
$ws.Range("E32").Value = "Pending"
$ws.Range("F32").Value = 14.21
$ws.Range("G32").Value = 1.42
$ws.Range("H32").Value = "SAHS"

$ws.Range("A33").Value = 45933.590578703705
$ws.Range("B33").Value = "M.15608.6485.47628"
$ws.Range("C33").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D33").Value = "Sale (iOS App)"
$ws.Range("E33").Value = "Pending"
$ws.Range("F33").Value = 12.63
$ws.Range("G33").Value = 1.26
$ws.Range("H33").Value = "SAHS"

$ws.Range("A34").Value = 45933.571388888886
$ws.Range("B34").Value = "M.15608.6485.37454"
$ws.Range("C34").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D34").Value = "Sale (iOS App)"
$ws.Range("E34").Value = "Pending"
$ws.Range("F34").Value = 36.3
$ws.Range("G34").Value = 3.63
$ws.Range("H34").Value = "SAHS"

$ws.Range("A35").Value = 45933.561574074076
$ws.Range("B35").Value = "M.15608.6485.32746"
$ws.Range("C35").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D35").Value = "Sale (iOS App)"
$ws.Range("E35").Value = "Pending"
$ws.Range("F35").Value = 110.49
$ws.Range("G35").Value = 11.05
$ws.Range("H35").Value = "SAHS"

$ws.Range("A36").Value = 45933.543275462966
$ws.Range("B36").Value = "M.15608.6485.23786"
$ws.Range("C36").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D36").Value = "Sale (iOS App)"
$ws.Range("E36").Value = "Pending"
$ws.Range("F36").Value = 39.46
$ws.Range("G36").Value = 3.95
$ws.Range("H36").Value = "SAHS"

$ws.Range("A37").Value = 45933.40354166667
$ws.Range("B37").Value = "M.15608.6484.1226209"
$ws.Range("C37").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D37").Value = "Sale (iOS App)"
$ws.Range("E37").Value = "Pending"
$ws.Range("F37").Value = 58.4
$ws.Range("G37").Value = 5.84
$ws.Range("H37").Value = "SAHS"

$ws.Range("A38").Value = 45933.2671412037
$ws.Range("B38").Value = "M.15608.6484.1092631"
$ws.Range("C38").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D38").Value = "Sale (iOS App)"
$ws.Range("E38").Value = "Pending"
$ws.Range("F38").Value = 64.72
$ws.Range("G38").Value = 6.47
$ws.Range("H38").Value = "SAHS"

$ws.Range("A39").Value = 45932.96487268519
$ws.Range("B39").Value = "M.15608.6484.540290"
$ws.Range("C39").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D39").Value = "Sale (iOS App)"
$ws.Range("E39").Value = "Pending"
$ws.Range("F39").Value = 110.49
$ws.Range("G39").Value = 11.05
$ws.Range("H39").Value = "SAHS"

$ws.Range("A40").Value = 45932.92831018518
$ws.Range("B40").Value = "M.15608.6484.452004"
$ws.Range("C40").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D40").Value = "Sale (iOS App)"
$ws.Range("E40").Value = "Pending"
$ws.Range("F40").Value = 29.63
$ws.Range("G40").Value = 2.96
$ws.Range("H40").Value = "D35"

$ws.Range("A41").Value = 45932.81166666667
$ws.Range("B41").Value = "M.15608.6484.256134"
$ws.Range("C41").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D41").Value = "Sale (iOS App)"
$ws.Range("E41").Value = "Pending"
$ws.Range("F41").Value = 122.24
$ws.Range("G41").Value = 12.22
$ws.Range("H41").Value = "SAHS15"

$ws.Range("A42").Value = 45932.68298611111
$ws.Range("B42").Value = "M.15608.6484.104465"
$ws.Range("C42").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D42").Value = "Sale (iOS App)"
$ws.Range("E42").Value = "Pending"
$ws.Range("F42").Value = 36.3
$ws.Range("G42").Value = 3.63
$ws.Range("H42").Value = "SAHS"

$ws.Range("A43").Value = 45932.6625462963
$ws.Range("B43").Value = "M.15608.6484.88865"
$ws.Range("C43").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D43").Value = "Sale (iOS App)"
$ws.Range("E43").Value = "Pending"
$ws.Range("F43").Value = 45.78
$ws.Range("G43").Value = 4.58
$ws.Range("H43").Value = "SAHS"

$ws.Range("A44").Value = 45932.636030092595
$ws.Range("B44").Value = "M.15608.6484.70793"
$ws.Range("C44").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D44").Value = "Sale (iOS App)"
$ws.Range("E44").Value = "Pending"
$ws.Range("F44").Value = 45.78
$ws.Range("G44").Value = 4.58
$ws.Range("H44").Value = "SAHS"

$ws.Range("A45").Value = 45932.540613425925
$ws.Range("B45").Value = "M.15608.6484.22125"
$ws.Range("C45").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D45").Value = "Sale (iOS App)"
$ws.Range("E45").Value = "Pending"
$ws.Range("F45").Value = 52.09
$ws.Range("G45").Value = 5.21
$ws.Range("H45").Value = "SAHS"

$ws.Range("A46").Value = 45932.51752314815
$ws.Range("B46").Value = "M.15608.6484.10213"
$ws.Range("C46").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D46").Value = "Sale (iOS App)"
$ws.Range("E46").Value = "Pending"
$ws.Range("F46").Value = 28.41
$ws.Range("G46").Value = 2.84
$ws.Range("H46").Value = "SAHS"

$ws.Range("A47").Value = 45932.408530092594
$ws.Range("B47").Value = "M.15608.6483.1228977"
$ws.Range("C47").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D47").Value = "Sale (iOS App)"
$ws.Range("E47").Value = "Pending"
$ws.Range("F47").Value = 61.1
$ws.Range("G47").Value = 6.11
$ws.Range("H47").Value = "SAHS15"

$ws.Range("A48").Value = 45932.39805555555
$ws.Range("B48").Value = "M.15608.6483.1212476"
$ws.Range("C48").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D48").Value = "Sale (iOS App)"
$ws.Range("E48").Value = "Pending"
$ws.Range("F48").Value = 17.36
$ws.Range("G48").Value = 1.74
$ws.Range("H48").Value = "SAHS"

$ws.Range("A49").Value = 45932.14304398148
$ws.Range("B49").Value = "M.15608.6483.967351"
$ws.Range("C49").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D49").Value = "Sale (iOS App)"
$ws.Range("E49").Value = "Pending"
$ws.Range("F49").Value = 59.49
$ws.Range("G49").Value = 5.95
$ws.Range("H49").Value = "SAHS15"

$ws.Range("A50").Value = 45931.955972222226
$ws.Range("B50").Value = "M.15608.6483.626399"
$ws.Range("C50").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D50").Value = "Sale (iOS App)"
$ws.Range("E50").Value = "Pending"
$ws.Range("F50").Value = 14.21
$ws.Range("G50").Value = 1.42
$ws.Range("H50").Value = "SAHS"

$ws.Range("A51").Value = 45931.89545138889
$ws.Range("B51").Value = "M.15608.6483.512475"
$ws.Range("C51").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D51").Value = "Sale (iOS App)"
$ws.Range("E51").Value = "Pending"
$ws.Range("F51").Value = 55.25
$ws.Range("G51").Value = 5.52
$ws.Range("H51").Value = "SAHS"

$ws.Range("A52").Value = 45931.81099537037
$ws.Range("B52").Value = "M.15608.6483.352117"
$ws.Range("C52").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D52").Value = "Sale (iOS App)"
$ws.Range("E52").Value = "Pending"
$ws.Range("F52").Value = 28.41
$ws.Range("G52").Value = 2.84
$ws.Range("H52").Value = "SAHS"

$ws.Range("A53").Value = 45931.80322916667
$ws.Range("B53").Value = "M.15608.6483.339235"
$ws.Range("C53").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D53").Value = "Sale (iOS App)"
$ws.Range("E53").Value = "Pending"
$ws.Range("F53").Value = 28.41
$ws.Range("G53").Value = 2.84
$ws.Range("H53").Value = "SAHS"

$ws.Range("A54").Value = 45931.73730324074
$ws.Range("B54").Value = "M.15608.6483.203773"
$ws.Range("C54").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D54").Value = "Sale (iOS App)"
$ws.Range("E54").Value = "Pending"
$ws.Range("F54").Value = 28.41
$ws.Range("G54").Value = 2.84
$ws.Range("H54").Value = "SAHS"

$ws.Range("A55").Value = 45931.677453703705
$ws.Range("B55").Value = "M.15608.6483.112339"
$ws.Range("C55").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D55").Value = "Sale (iOS App)"
$ws.Range("E55").Value = "Pending"
$ws.Range("F55").Value = 14.21
$ws.Range("G55").Value = 1.42
$ws.Range("H55").Value = "SAHS"

$ws.Range("A56").Value = 45931.67476851852
$ws.Range("B56").Value = "M.15608.6483.110261"
$ws.Range("C56").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D56").Value = "Sale (iOS App)"
$ws.Range("E56").Value = "Pending"
$ws.Range("F56").Value = 14.21
$ws.Range("G56").Value = 1.42
$ws.Range("H56").Value = "SAHS"

$ws.Range("A57").Value = 45931.63186342592
$ws.Range("B57").Value = "M.15608.6483.79557"
$ws.Range("C57").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D57").Value = "Sale (iOS App)"
$ws.Range("E57").Value = "Pending"
$ws.Range("F57").Value = 44.2
$ws.Range("G57").Value = 4.42
$ws.Range("H57").Value = "SAHS"

$ws.Range("A58").Value = 45931.50315972222
$ws.Range("B58").Value = "M.15608.6483.1918"
$ws.Range("C58").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D58").Value = "Sale (iOS App)"
$ws.Range("E58").Value = "Pending"
$ws.Range("F58").Value = 55.25
$ws.Range("G58").Value = 5.52
$ws.Range("H58").Value = "SAHS"

$ws.Range("A59").Value = 45931.48385416667
$ws.Range("B59").Value = "M.15608.6482.1391354"
$ws.Range("C59").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D59").Value = "Sale (iOS App)"
$ws.Range("E59").Value = "Pending"
$ws.Range("F59").Value = 44.2
$ws.Range("G59").Value = 4.42
$ws.Range("H59").Value = "SAHS"

$ws.Range("A60").Value = 45931.46702546296
$ws.Range("B60").Value = "M.15608.6482.1375918"
$ws.Range("C60").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D60").Value = "Sale (iOS App)"
$ws.Range("E60").Value = "Pending"
$ws.Range("F60").Value = 14.21
$ws.Range("G60").Value = 1.42
$ws.Range("H60").Value = "SAHS"

$ws.Range("A61").Value = 45931.36540509259
$ws.Range("B61").Value = "M.15608.6482.1256968"
$ws.Range("C61").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D61").Value = "Sale (iOS App)"
$ws.Range("E61").Value = "Pending"
$ws.Range("F61").Value = 39.46
$ws.Range("G61").Value = 3.95
$ws.Range("H61").Value = "SAHS"

$ws.Range("A62").Value = 45931.341099537036
$ws.Range("B62").Value = "M.15608.6482.1239968"
$ws.Range("C62").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D62").Value = "Sale (iOS App)"
$ws.Range("E62").Value = "Pending"
$ws.Range("F62").Value = 41.04
$ws.Range("G62").Value = 4.1
$ws.Range("H62").Value = "SAHS"

$ws.Range("A63").Value = 45931.31701388889
$ws.Range("B63").Value = "M.15608.6482.1209917"
$ws.Range("C63").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D63").Value = "Sale (iOS App)"
$ws.Range("E63").Value = "Pending"
$ws.Range("F63").Value = 110.49
$ws.Range("G63").Value = 11.05
$ws.Range("H63").Value = "SAHS"

$ws.Range("A64").Value = 45931.30383101852
$ws.Range("B64").Value = "M.15608.6482.1198944"
$ws.Range("C64").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D64").Value = "Sale (iOS App)"
$ws.Range("E64").Value = "Pending"
$ws.Range("F64").Value = 17.89
$ws.Range("G64").Value = 1.79
$ws.Range("H64").Value = "SAHS15"

$ws.Range("A65").Value = 45931.25746527778
$ws.Range("B65").Value = "M.15608.6482.1145343"
$ws.Range("C65").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D65").Value = "Sale (iOS App)"
$ws.Range("E65").Value = "Pending"
$ws.Range("F65").Value = 146.1
$ws.Range("G65").Value = 14.61
$ws.Range("H65").Value = "SAHS15"

$ws.Range("A66").Value = 45930.92653935185
$ws.Range("B66").Value = "M.15608.6482.562919"
$ws.Range("C66").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D66").Value = "Sale (iOS App)"
$ws.Range("E66").Value = "Pending"
$ws.Range("F66").Value = 55.24
$ws.Range("G66").Value = 5.52
$ws.Range("H66").Value = "SAHS"

$ws.Range("A67").Value = 45930.74658564815
$ws.Range("B67").Value = "M.15608.6482.211929"
$ws.Range("C67").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D67").Value = "Sale (Android App)"
$ws.Range("E67").Value = "Pending"
$ws.Range("F67").Value = 51.83
$ws.Range("G67").Value = 5.18
$ws.Range("H67").Value = "SAHS"

$ws.Range("A68").Value = 45930.454872685186
$ws.Range("B68").Value = "M.15608.6481.1355187"
$ws.Range("C68").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D68").Value = "Sale (iOS App)"
$ws.Range("E68").Value = "Pending"
$ws.Range("F68").Value = 265.32
$ws.Range("G68").Value = 26.53
$ws.Range("H68").Value = "SAHS15"

$ws.Range("A69").Value = 45930.279074074075
$ws.Range("B69").Value = "15608.6481.1176060"
$ws.Range("C69").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D69").Value = "Sale"
$ws.Range("E69").Value = "Pending"
$ws.Range("F69").Value = 36.3
$ws.Range("G69").Value = 3.63
$ws.Range("H69").Value = "D26"

$ws.Range("A70").Value = 45930.22298611111
$ws.Range("B70").Value = "M.15608.6481.1099989"
$ws.Range("C70").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D70").Value = "Sale (iOS App)"
$ws.Range("E70").Value = "Pending"
$ws.Range("F70").Value = 25.25
$ws.Range("G70").Value = 2.53
$ws.Range("H70").Value = "SAHS"

$ws.Range("A71").Value = 45930.125601851854
$ws.Range("B71").Value = "M.15608.6481.957577"
$ws.Range("C71").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D71").Value = "Sale (Android App)"
$ws.Range("E71").Value = "Pending"
$ws.Range("F71").Value = 28.02
$ws.Range("G71").Value = 2.8
$ws.Range("H71").Value = "SAHS15"

$ws.Range("A72").Value = 45930.09780092593
$ws.Range("B72").Value = "M.15608.6481.892957"
$ws.Range("C72").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D72").Value = "Sale (iOS App)"
$ws.Range("E72").Value = "Pending"
$ws.Range("F72").Value = 58.11
$ws.Range("G72").Value = 5.81
$ws.Range("H72").Value = "SAHS15"

$ws.Range("A73").Value = 45930.095509259256
$ws.Range("B73").Value = "M.15608.6481.891385"
$ws.Range("C73").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D73").Value = "Sale (iOS App)"
$ws.Range("E73").Value = "Pending"
$ws.Range("F73").Value = 104.34
$ws.Range("G73").Value = 10.43
$ws.Range("H73").Value = "SAHS15"

$ws.Range("A74").Value = 45930.083402777775
$ws.Range("B74").Value = "M.15608.6481.874839"
$ws.Range("C74").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D74").Value = "Sale (iOS App)"
$ws.Range("E74").Value = "Pending"
$ws.Range("F74").Value = 110.48
$ws.Range("G74").Value = 11.05
$ws.Range("H74").Value = "SAHS"

$ws.Range("A75").Value = 45930.03109953704
$ws.Range("B75").Value = "M.15608.6481.754754"
$ws.Range("C75").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D75").Value = "Sale (iOS App)"
$ws.Range("E75").Value = "Pending"
$ws.Range("F75").Value = 110.48
$ws.Range("G75").Value = 11.05
$ws.Range("H75").Value = "SAHS"

$ws.Range("A76").Value = 45929.79486111111
$ws.Range("B76").Value = "M.15608.6481.339089"
$ws.Range("C76").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D76").Value = "Sale (iOS App)"
$ws.Range("E76").Value = "Pending"
$ws.Range("F76").Value = 36.3
$ws.Range("G76").Value = 3.63
$ws.Range("H76").Value = "SAHS"

$ws.Range("A77").Value = 45929.791863425926
$ws.Range("B77").Value = "M.15608.6481.328255"
$ws.Range("C77").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D77").Value = "Sale (iOS App)"
$ws.Range("E77").Value = "Pending"
$ws.Range("F77").Value = 110.48
$ws.Range("G77").Value = 11.05
$ws.Range("H77").Value = "SAHS"

$ws.Range("A78").Value = 45929.6397337963
$ws.Range("B78").Value = "M.15608.6481.81653"
$ws.Range("C78").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D78").Value = "Sale (iOS App)"
$ws.Range("E78").Value = "Pending"
$ws.Range("F78").Value = 137.13
$ws.Range("G78").Value = 13.71
$ws.Range("H78").Value = "SAHS15"

$ws.Range("A79").Value = 45929.61986111111
$ws.Range("B79").Value = "M.15608.6481.70161"
$ws.Range("C79").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D79").Value = "Sale (iOS App)"
$ws.Range("E79").Value = "Pending"
$ws.Range("F79").Value = 52.08
$ws.Range("G79").Value = 5.21
$ws.Range("H79").Value = "SAHS"

$ws.Range("A80").Value = 45929.51945601852
$ws.Range("B80").Value = "M.15608.6481.12232"
$ws.Range("C80").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D80").Value = "Sale (iOS App)"
$ws.Range("E80").Value = "Pending"
$ws.Range("F80").Value = 110.48
$ws.Range("G80").Value = 11.05
$ws.Range("H80").Value = "SAHS"

$ws.Range("A81").Value = 45929.51195601852
$ws.Range("B81").Value = "M.15608.6481.9150"
$ws.Range("C81").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D81").Value = "Sale (iOS App)"
$ws.Range("E81").Value = "Pending"
$ws.Range("F81").Value = 14.2
$ws.Range("G81").Value = 1.42
$ws.Range("H81").Value = "SAHS"

$ws.Range("A82").Value = 45929.06377314815
$ws.Range("B82").Value = "M.15608.6480.786055"
$ws.Range("C82").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D82").Value = "Sale (iOS App)"
$ws.Range("E82").Value = "Pending"
$ws.Range("F82").Value = 52.15
$ws.Range("G82").Value = 5.22
$ws.Range("H82").Value = "SAHS15"

$ws.Range("A83").Value = 45929.05862268519
$ws.Range("B83").Value = "M.15608.6480.778925"
$ws.Range("C83").Value = "Airalo - The World's First eSIM Store"
$ws.Range("D83").Value = "Sale (iOS App)"
$ws.Range("E83").Value = "Pending"
$ws.Range("F83").Value = 56.82
$ws.Range("G83").Value = 5.68
$ws.Range("H83").Value = "SAHS"

# Step 3: update the summary row (84) totals (other cells in the row stay blank/unchanged).
$ws.Range("F84").Value = 4740.73
$ws.Range("G84").Value = 474.05
